# Replace the old "name/file/alt/calorie/gluten/lactose" header table
# (spread over columns A:G) with a new two-column cake/calorie table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out everything that is currently used (A1:G9) so no stray data
# remains from the old layout.
$ws.Range("A1:G9").Clear()

# New data: cake name in column A, calorie value in column B.
$data = @(
    @("Áfonya torta", "585 Kcal"),
    @("Csoki torta", "394 Kcal"),
    @("Fánk", "545 Kcal"),
    @("Krémes mocsi", "221 Kcal"),
    @("Kuglóf", "438 Kcal"),
    @("Macaron", "990 Kcal"),
    @("Muffin", "550 Kcal"),
    @("Piskota", "476 Kcal")
)

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Last two rows were typed as both names first, then both calorie values
# (matches the shared-string insertion order recorded in the workbook).
$ws.Cells.Item(7, 1).Value = $data[6][0]
$ws.Cells.Item(8, 1).Value = $data[7][0]
$ws.Cells.Item(7, 2).Value = $data[6][1]
$ws.Cells.Item(8, 2).Value = $data[7][1]

# Column widths matching the new layout (closest achievable values to the
# saved 14.85546875 / 11.7109375 character widths).
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 10.86

# Selection moves to B8, matching the saved view state.
$ws.Range("B8").Select()
